$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @("FAPs","Fgf15","Fgfr4","ECs",1,0.3333333333333333,0.029478,0.088434,0.1535387136874709,0.1535387136874709,3,1,0.8372556666666666,2.511767,0.04323017762831461,0.04323017762831462,0.024680622542,0.222125602878,0.006637505865532307,0.006637505865532308),
    @("FAPs","Fgf15","Fgfr4","FAPs",1,0.3333333333333333,0.029478,0.088434,0.1535387136874709,0.1535387136874709,2,0.6666666666666666,0.110372,0.331116,0.005698858013333649,0.00569885801333365,0.003253545816,0.029281912344,0.0008749953288547844,0.0008749953288547846),
    @("FAPs","Fgf15","Fgfr4","M2",1,0.3333333333333333,0.029478,0.088434,0.1535387136874709,0.1535387136874709,1,0.3333333333333333,0.01411566666666667,0.042347,0.0007288368435552497,0.0007288368435552499,0.000416101622,0.003744914598,0.0001119046714475095,0.0001119046714475095),
    @("FAPs","Fgf15","Fgfr4","sCs",1,0.3333333333333333,0.029478,0.088434,0.1535387136874709,0.1535387136874709,3,1,18.40564566666667,55.216937,0.9503421275147964,0.9503421275147965,0.5425616229620001,4.883054606658,0.1459143078216363,0.1459143078216363),
    @("sCs","Fgf15","Fgfr4","ECs",3,1,0.1625126666666667,0.487538,0.846461286312529,0.846461286312529,3,1,0.8372556666666666,2.511767,0.04323017762831461,0.04323017762831462,0.1360646510717778,1.224581859646,0.0365926717627823,0.03659267176278231),
    @("sCs","Fgf15","Fgfr4","FAPs",3,1,0.1625126666666667,0.487538,0.846461286312529,0.846461286312529,2,0.6666666666666666,0.110372,0.331116,0.005698858013333649,0.00569885801333365,0.01793684804533333,0.161431632408,0.004823862684478864,0.004823862684478866),
    @("sCs","Fgf15","Fgfr4","M2",3,1,0.1625126666666667,0.487538,0.846461286312529,0.846461286312529,1,0.3333333333333333,0.01411566666666667,0.042347,0.0007288368435552497,0.0007288368435552499,0.002293974631777778,0.020645771686,0.0006169321721077401,0.0006169321721077403),
    @("sCs","Fgf15","Fgfr4","sCs",3,1,0.1625126666666667,0.487538,0.846461286312529,0.846461286312529,3,1,18.40564566666667,55.216937,0.9503421275147964,0.9503421275147965,2.991150559011778,26.920355031106,0.80442781969316,0.8044278196931601),
)

for ($i = 0; $i -lt $rowsData.Length; $i++) {
    $r = $i + 2
    $rowVals = $rowsData[$i]
    for ($c = 1; $c -le $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}

Write-Host "Updated rows 2-9"
